# Update countries & provincias Spain
# Applies the 20-Abril-2020 01:22 data refresh to the "Pais" sheet:
#   - Nigeria's case counts updated, re-sorted just below Niger (was below Kirguistan)
#   - Libia's case counts updated, re-sorted just below Maldivas (was below Guinea-Bisau)
#   - Santo Tome y Principe re-sorted just above Sudan del Sur (tied totals)
#   - Estados Unidos (row 4) totals refreshed
#   - Footer timestamp advanced from 00:52 to 01:22

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4: Estados Unidos - refreshed totals ---------------------------
$ws.Range("B4").Value = 763579
$ws.Range("C4").Value = 24787
$ws.Range("D4").Value = 70938
$ws.Range("E4").Value = 652117
$ws.Range("F4").Value = 13566
$ws.Range("G4").Value = 1510
$ws.Range("H4").Value = 40524

# --- Rows 96-100: Nigeria jumps above Guinea/Burkina Faso/Albania/Kirguistan,
#     which each shift down one place keeping their previous figures --------
$ws.Range("A96").Value = "Nigeria"
$ws.Range("B96").Value = 627
$ws.Range("C96").Value = 85
$ws.Range("D96").Value = 170
$ws.Range("E96").Value = 436
$ws.Range("F96").Value = 2
$ws.Range("G96").Value = 2
$ws.Range("H96").Value = 21

$ws.Range("A97").Value = "Guinea"
$ws.Range("B97").Value = 579
$ws.Range("C97").Value = 61
$ws.Range("D97").Value = 87
$ws.Range("E97").Value = 487
$ws.Range("F97").Value = 0
$ws.Range("G97").Value = 2
$ws.Range("H97").Value = 5

$ws.Range("A98").Value = "Burkina Faso"
$ws.Range("B98").Value = 576
$ws.Range("C98").Value = 11
$ws.Range("D98").Value = 338
$ws.Range("E98").Value = 202
$ws.Range("F98").Value = 0
$ws.Range("G98").Value = 0
$ws.Range("H98").Value = 36

$ws.Range("A99").Value = "Albania"
$ws.Range("B99").Value = 562
$ws.Range("C99").Value = 14
$ws.Range("D99").Value = 314
$ws.Range("E99").Value = 222
$ws.Range("F99").Value = 5
$ws.Range("G99").Value = 0
$ws.Range("H99").Value = 26

$ws.Range("A100").Value = "Kirguistan"
$ws.Range("B100").Value = 554
$ws.Range("C100").Value = 48
$ws.Range("D100").Value = 133
$ws.Range("E100").Value = 416
$ws.Range("F100").Value = 5
$ws.Range("G100").Value = 0
$ws.Range("H100").Value = 5

# --- Rows 160-161: Libia jumps above Guinea-Bisau, which shifts down one
#     place keeping its previous figures ------------------------------------
$ws.Range("A160").Value = "Libia"
$ws.Range("B160").Value = 51
$ws.Range("C160").Value = 2
$ws.Range("D160").Value = 11
$ws.Range("E160").Value = 39
$ws.Range("F160").Value = 0
$ws.Range("G160").Value = 0
$ws.Range("H160").Value = 1

$ws.Range("A161").Value = "Guinea-Bisau"
$ws.Range("B161").Value = 50
$ws.Range("C161").Value = 4
$ws.Range("D161").Value = 3
$ws.Range("E161").Value = 47
$ws.Range("F161").Value = 0
$ws.Range("G161").Value = 0
$ws.Range("H161").Value = 0

# --- Rows 211-212: Santo Tome y Principe swaps ahead of Sudan del Sur
#     (both tied on totals, figures unchanged) ------------------------------
$ws.Range("A211").Value = "Santo Tome y Principe"
$ws.Range("A212").Value = "Sudan del Sur"

# --- Footer timestamp ------------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 20 de Abril de 2020 a las 01:22"
